# Add team record (Wins / Losses / Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new headers in AD1:AF1, copying the existing header
# style (bold, bordered, centered) from an existing header cell like A1.
$ws.Range("A1").Copy($ws.Range("AD1"))
$ws.Range("A1").Copy($ws.Range("AE1"))
$ws.Range("A1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-58: every row gets the same team record (77 wins, 85 losses, 0 ties).
$lastRow = 58
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 77
    $ws.Cells.Item($r, 31).Value = 85
    $ws.Cells.Item($r, 32).Value = 0
}
